$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextValue 2 4 "304.33"
Set-TextValue 2 5 "1.07%"
Set-TextValue 2 7 "22"

Set-TextValue 3 4 "35.56"
Set-TextValue 3 5 "10.55%"
Set-TextValue 3 7 "22"

Set-TextValue 4 4 "5.080"
Set-TextValue 4 5 "1.75%"
Set-TextValue 4 7 "22"

Set-TextValue 5 4 "0.07818"
Set-TextValue 5 5 "0.96%"
Set-TextValue 5 7 "22"

Set-TextValue 6 4 "2.265"
Set-TextValue 6 5 "-0.87%"
Set-TextValue 6 7 "22"

Set-TextValue 7 5 "1.59%"
Set-TextValue 7 7 "22"

Set-TextValue 8 4 "4.044"
Set-TextValue 8 5 "5.52%"
Set-TextValue 8 7 "22"

Set-TextValue 9 4 "0.9291"
Set-TextValue 9 5 "-0.18%"
Set-TextValue 9 7 "22"

Set-TextValue 10 4 "0.09306"
Set-TextValue 10 5 "-7.98%"
Set-TextValue 10 7 "22"

Set-TextValue 11 4 "0.1839"
Set-TextValue 11 5 "3.92%"
Set-TextValue 11 7 "22"

Set-TextValue 12 4 "0.08606"
Set-TextValue 12 5 "1.64%"
Set-TextValue 12 7 "22"

Set-TextValue 13 4 "0.03768"
Set-TextValue 13 5 "15.12%"
Set-TextValue 13 7 "22"

Set-TextValue 14 4 "0.09967"
Set-TextValue 14 5 "0.68%"
Set-TextValue 14 7 "22"

Set-TextValue 15 4 "0.001486"
Set-TextValue 15 5 "-0.41%"
Set-TextValue 15 7 "22"

Set-TextValue 16 4 "0.005640"
Set-TextValue 16 5 "-0.72%"
Set-TextValue 16 7 "22"

Set-TextValue 17 5 "-1.12%"
Set-TextValue 17 7 "22"

Set-TextValue 18 5 "-5.82%"
Set-TextValue 18 7 "22"

Set-TextValue 19 5 "1.18%"
Set-TextValue 19 7 "22"

Set-TextValue 20 4 "0.1321"
Set-TextValue 20 5 "-1.35%"
Set-TextValue 20 7 "22"

Set-TextValue 21 4 "4.560"
Set-TextValue 21 5 "4.69%"
Set-TextValue 21 7 "22"

Set-TextValue 22 4 "0.2236"
Set-TextValue 22 5 "6.92%"
Set-TextValue 22 7 "22"

Set-TextValue 23 4 "0.04678"
Set-TextValue 23 5 "2.39%"
Set-TextValue 23 7 "22"

Set-TextValue 24 5 "1.16%"
Set-TextValue 24 7 "22"

Set-TextValue 25 4 "0.004548"
Set-TextValue 25 5 "4.09%"
Set-TextValue 25 7 "22"

Set-TextValue 26 5 "0.66%"
Set-TextValue 26 7 "22"

Set-TextValue 27 5 "-20.12%"
Set-TextValue 27 7 "22"

Set-TextValue 28 7 "22"

Set-TextValue 29 7 "22"

Set-TextValue 30 7 "22"

Set-TextValue 31 7 "22"

Set-TextValue 32 7 "22"

Set-TextValue 33 7 "22"

Set-TextValue 34 7 "22"

Set-TextValue 35 7 "22"

Set-TextValue 36 7 "22"

Set-TextValue 37 7 "22"

Set-TextValue 38 7 "22"

Set-TextValue 39 4 "0.01785"
Set-TextValue 39 5 "5.15%"
Set-TextValue 39 7 "22"

Set-TextValue 40 4 "0.04722"
Set-TextValue 40 7 "22"

Set-TextValue 41 4 "0.007889"
Set-TextValue 41 5 "2.20%"
Set-TextValue 41 7 "22"

Set-TextValue 42 5 "2.08%"
Set-TextValue 42 7 "22"

Set-TextValue 43 4 "0.007990"
Set-TextValue 43 5 "-18.37%"
Set-TextValue 43 7 "22"

Set-TextValue 44 5 "8.61%"
Set-TextValue 44 7 "22"

Set-TextValue 45 4 "0.009067"
Set-TextValue 45 5 "-6.16%"
Set-TextValue 45 7 "22"

Set-TextValue 46 4 "0.00006193"
Set-TextValue 46 5 "1.94%"
Set-TextValue 46 7 "22"

Set-TextValue 47 4 "0.00000000750"
Set-TextValue 47 5 "0.61%"
Set-TextValue 47 7 "22"

Set-TextValue 48 4 "4.086"
Set-TextValue 48 5 "53.95%"
Set-TextValue 48 7 "22"

Set-TextValue 49 4 "0.002690"
Set-TextValue 49 5 "35.09%"
Set-TextValue 49 7 "22"

Set-TextValue 50 4 "0.00002100"
Set-TextValue 50 5 "0.61%"
Set-TextValue 50 7 "22"

Set-TextValue 51 4 "0.0002000"
Set-TextValue 51 5 "0.61%"
Set-TextValue 51 7 "22"
